$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (rows 5, 6, 7) mirroring the existing rows
$data = @(
    @(42602.514305555553, "Noun", 4901, 1653, 289, 38, 15, 71, 28, 0, 1, 0, 100),
    @(42602.516250000001, "Noun", 5675, 1411, 237, 33, 14, 69, 29, 0, 1, 0, 100),
    @(42602.517106481479, "Noun", 5290, 1654, 289, 38, 15, 71, 28, 0, 1, 0, 100)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

# Copy the style (number format) from an existing date cell (A2) down to the
# new date cells so we reuse the existing style entry instead of creating a
# new one.
$ws.Range("A2").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
